# 自动更新Excel文件 - 2025-12-28 23:13:07
# Recalculates column E ("剩余" / days remaining) for every data row based on
# column D ("总天" / total days) and column F ("开始时间" / start date, stored
# as an integer YYYYMMDD), relative to the reference date 2025-12-29.
#
# Rule (mirrors the automated nightly recompute job):
#   elapsed = referenceDate - startDate   (in days)
#   if elapsed < 0            -> row not due yet, leave untouched
#   if elapsed >= totalDays    -> the cycle rolled over: reset start date to
#                                 the reference date and remaining = totalDays
#   else                       -> remaining = totalDays - elapsed  (F unchanged)
#
# Rows whose start date cannot be parsed as an 8-digit YYYYMMDD value
# (a data-entry error) are left untouched, exactly as the source workbook
# has them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$refDate = Get-Date -Year 2025 -Month 12 -Day 29 -Hour 0 -Minute 0 -Second 0
$refOA = $refDate.ToOADate()
$refSerial = 20251229

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {

    $totalDaysVal = $ws.Cells.Item($r, 4).Value2   # column D
    $startVal     = $ws.Cells.Item($r, 6).Value2   # column F

    if ($null -eq $totalDaysVal -or $null -eq $startVal) {
        continue
    }

    $totalDays = [int]$totalDaysVal
    $startStr = [string]([int64]$startVal)

    if ($startStr.Length -ne 8) {
        # malformed date (e.g. "202510929") - skip, matches source data
        continue
    }

    $year  = [int]$startStr.Substring(0,4)
    $month = [int]$startStr.Substring(4,2)
    $day   = [int]$startStr.Substring(6,2)

    if ($month -lt 1 -or $month -gt 12 -or $day -lt 1 -or $day -gt 31) {
        continue
    }

    $startDate = Get-Date -Year $year -Month $month -Day $day -Hour 0 -Minute 0 -Second 0
    $startOA = $startDate.ToOADate()

    $elapsed = [int]([math]::Round($refOA - $startOA))

    if ($elapsed -lt 0) {
        # start date still in the future relative to the reference date
        continue
    }

    if ($elapsed -ge $totalDays) {
        # cycle complete - restart the clock from the reference date
        $ws.Cells.Item($r, 5).Value = $totalDays
        $ws.Cells.Item($r, 6).Value = $refSerial
    } else {
        $ws.Cells.Item($r, 5).Value = ($totalDays - $elapsed)
    }
}
